$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 3).Value = 1.029309600808322
$ws.Cells.Item(2, 4).Value = 1.045743174522805
$ws.Cells.Item(2, 5).Value = 1.04082740150638
$ws.Cells.Item(2, 6).Value = 1.049676014374546
$ws.Cells.Item(2, 9).Value = 1.061202169094939
$ws.Cells.Item(2, 10).Value = 1.050647662194486
$ws.Cells.Item(2, 11).Value = 1.056599999229192
$ws.Cells.Item(2, 12).Value = 1.051745804848505
$ws.Cells.Item(2, 13).Value = 1.060484090360038
$ws.Cells.Item(2, 14).Value = 1.020275197344785
$ws.Cells.Item(2, 15).Value = 1.03
$ws.Cells.Item(2, 16).Value = 1.056441318569447
$ws.Cells.Item(2, 18).Value = 1.051091251677588

# Row 3
$ws.Cells.Item(3, 3).Value = 1.033202069142158
$ws.Cells.Item(3, 4).Value = 1.048419679621628
$ws.Cells.Item(3, 5).Value = 1.043864331240119
$ws.Cells.Item(3, 6).Value = 1.052527855518423
$ws.Cells.Item(3, 9).Value = 1.062375959116271
$ws.Cells.Item(3, 10).Value = 1.052826789526942
$ws.Cells.Item(3, 11).Value = 1.058468353078088
$ws.Cells.Item(3, 12).Value = 1.053965113914912
$ws.Cells.Item(3, 13).Value = 1.062530011640093
$ws.Cells.Item(3, 14).Value = 1.021038646628351
$ws.Cells.Item(3, 15).Value = 1.03
$ws.Cells.Item(3, 16).Value = 1.058060491577847
$ws.Cells.Item(3, 18).Value = 1.052409677063856

# Row 4
$ws.Cells.Item(4, 3).Value = 1.035676648367626
$ws.Cells.Item(4, 4).Value = 1.050125024824866
$ws.Cells.Item(4, 5).Value = 1.045800789304435
$ws.Cells.Item(4, 6).Value = 1.054347935904692
$ws.Cells.Item(4, 9).Value = 1.063113086457444
$ws.Cells.Item(4, 10).Value = 1.054209766810827
$ws.Cells.Item(4, 11).Value = 1.059653904106261
$ws.Cells.Item(4, 12).Value = 1.055376162541683
$ws.Cells.Item(4, 13).Value = 1.063831861084189
$ws.Cells.Item(4, 14).Value = 1.021523085723655
$ws.Cells.Item(4, 15).Value = 1.03
$ws.Cells.Item(4, 16).Value = 1.059090793116107
$ws.Cells.Item(4, 18).Value = 1.053248815632907

# Row 5
$ws.Cells.Item(5, 3).Value = 1.036710892566961
$ws.Cells.Item(5, 4).Value = 1.050840562325938
$ws.Cells.Item(5, 5).Value = 1.046612190147232
$ws.Cells.Item(5, 6).Value = 1.055110866423504
$ws.Cells.Item(5, 9).Value = 1.063420572001756
$ws.Cells.Item(5, 10).Value = 1.054788931574238
$ws.Cells.Item(5, 11).Value = 1.060151586881193
$ws.Cells.Item(5, 12).Value = 1.055967472212649
$ws.Cells.Item(5, 13).Value = 1.064377633832069
$ws.Cells.Item(5, 14).Value = 1.021726836936376
$ws.Cells.Item(5, 15).Value = 1.03
$ws.Cells.Item(5, 16).Value = 1.059522724758019
$ws.Cells.Item(5, 18).Value = 1.053607748901822

# Row 6
$ws.Cells.Item(6, 3).Value = 1.03688904042655
$ws.Cells.Item(6, 4).Value = 1.050966186793695
$ws.Cells.Item(6, 5).Value = 1.046752891056045
$ws.Cells.Item(6, 6).Value = 1.055243066467535
$ws.Cells.Item(6, 9).Value = 1.063475379872486
$ws.Cells.Item(6, 10).Value = 1.054890781571659
$ws.Cells.Item(6, 11).Value = 1.060240613362897
$ws.Cells.Item(6, 12).Value = 1.056071216000172
$ws.Cells.Item(6, 13).Value = 1.064473369391328
$ws.Cells.Item(6, 14).Value = 1.021763760962669
$ws.Cells.Item(6, 15).Value = 1.03
$ws.Cells.Item(6, 16).Value = 1.059598491098235
$ws.Cells.Item(6, 18).Value = 1.05367928285801

# Row 7
$ws.Cells.Item(7, 3).Value = 1.035704223853917
$ws.Cells.Item(7, 4).Value = 1.050150317005738
$ws.Cells.Item(7, 5).Value = 1.045824745480966
$ws.Cells.Item(7, 6).Value = 1.054370114406539
$ws.Cells.Item(7, 9).Value = 1.063126593472147
$ws.Cells.Item(7, 10).Value = 1.054230927100636
$ws.Cells.Item(7, 11).Value = 1.059676114544491
$ws.Cells.Item(7, 12).Value = 1.055397033228931
$ws.Cells.Item(7, 13).Value = 1.06385100836039
$ws.Cells.Item(7, 14).Value = 1.02153351403898
$ws.Cells.Item(7, 15).Value = 1.03
$ws.Cells.Item(7, 16).Value = 1.059105946523299
$ws.Cells.Item(7, 18).Value = 1.053284272634615

# Row 8
$ws.Cells.Item(8, 3).Value = 1.030651708680276
$ws.Cells.Item(8, 4).Value = 1.046672925749468
$ws.Cells.Item(8, 5).Value = 1.041876347931048
$ws.Cells.Item(8, 6).Value = 1.050660161200289
$ws.Cells.Item(8, 9).Value = 1.061615280142348
$ws.Cells.Item(8, 10).Value = 1.051406651032532
$ws.Cells.Item(8, 11).Value = 1.057255725462307
$ws.Cells.Item(8, 12).Value = 1.052517450761204
$ws.Cells.Item(8, 13).Value = 1.061194998133698
$ws.Cells.Item(8, 14).Value = 1.02054504989777
$ws.Cells.Item(8, 15).Value = 1.03
$ws.Cells.Item(8, 16).Value = 1.057003942109353
$ws.Cells.Item(8, 18).Value = 1.051577408512006

# Row 9
$ws.Cells.Item(9, 3).Value = 1.021373671508313
$ws.Cells.Item(9, 4).Value = 1.040309198118861
$ws.Cells.Item(9, 5).Value = 1.034669267622255
$ws.Cells.Item(9, 6).Value = 1.043902230913669
$ws.Cells.Item(9, 9).Value = 1.058761109285724
$ws.Cells.Item(9, 10).Value = 1.046193607818038
$ws.Cells.Item(9, 11).Value = 1.052780468901844
$ws.Cells.Item(9, 12).Value = 1.047223832427577
$ws.Cells.Item(9, 13).Value = 1.056321032773578
$ws.Cells.Item(9, 14).Value = 1.018714981401561
$ws.Cells.Item(9, 15).Value = 1.03
$ws.Cells.Item(9, 16).Value = 1.053146601281269
$ws.Cells.Item(9, 18).Value = 1.048410077218292

# Row 10
$ws.Cells.Item(10, 3).Value = 1.015006508632126
$ws.Cells.Item(10, 4).Value = 1.035988109406395
$ws.Cells.Item(10, 5).Value = 1.029788576195711
$ws.Cells.Item(10, 6).Value = 1.039369390843533
$ws.Cells.Item(10, 9).Value = 1.056780491269738
$ws.Cells.Item(10, 10).Value = 1.042637130882439
$ws.Cells.Item(10, 11).Value = 1.049736629549426
$ws.Cells.Item(10, 12).Value = 1.043640041239432
$ws.Cells.Item(10, 13).Value = 1.053062429129004
$ws.Cells.Item(10, 14).Value = 1.017481059759038
$ws.Cells.Item(10, 15).Value = 1.03
$ws.Cells.Item(10, 16).Value = 1.050618366447801
$ws.Cells.Item(10, 18).Value = 1.046274511956011

# Row 11
$ws.Cells.Item(11, 3).Value = 1.012784565943969
$ws.Cells.Item(11, 4).Value = 1.034656153665791
$ws.Cells.Item(11, 5).Value = 1.028340717544906
$ws.Cells.Item(11, 6).Value = 1.038304600744824
$ws.Cells.Item(11, 9).Value = 1.056257334880455
$ws.Cells.Item(11, 10).Value = 1.041634996676806
$ws.Cells.Item(11, 11).Value = 1.048952819073953
$ws.Cells.Item(11, 12).Value = 1.042747254478299
$ws.Cells.Item(11, 13).Value = 1.052538580023374
$ws.Cells.Item(11, 14).Value = 1.017250735088883
$ws.Cells.Item(11, 15).Value = 1.03
$ws.Cells.Item(11, 16).Value = 1.050633978335286
$ws.Cells.Item(11, 18).Value = 1.045752711672792

# Row 12
$ws.Cells.Item(12, 3).Value = 1.012161900756003
$ws.Cells.Item(12, 4).Value = 1.034356369026454
$ws.Cells.Item(12, 5).Value = 1.028048707514999
$ws.Cells.Item(12, 6).Value = 1.038239095168286
$ws.Cells.Item(12, 9).Value = 1.056185927377153
$ws.Cells.Item(12, 10).Value = 1.041460319124462
$ws.Cells.Item(12, 11).Value = 1.048853683472655
$ws.Cells.Item(12, 12).Value = 1.042657570028649
$ws.Cells.Item(12, 13).Value = 1.052668617738648
$ws.Cells.Item(12, 14).Value = 1.017280116918191
$ws.Cells.Item(12, 15).Value = 1.03
$ws.Cells.Item(12, 16).Value = 1.051059120426481
$ws.Cells.Item(12, 18).Value = 1.045682621458397

# Row 13
$ws.Cells.Item(13, 3).Value = 1.012748716447473
$ws.Cells.Item(13, 4).Value = 1.034852914738849
$ws.Cells.Item(13, 5).Value = 1.028646313500698
$ws.Cells.Item(13, 6).Value = 1.038964070327251
$ws.Cells.Item(13, 9).Value = 1.056475276455466
$ws.Cells.Item(13, 10).Value = 1.041932473124717
$ws.Cells.Item(13, 11).Value = 1.049299624514151
$ws.Cells.Item(13, 12).Value = 1.043202268182285
$ws.Cells.Item(13, 13).Value = 1.053339347594822
$ws.Cells.Item(13, 14).Value = 1.017522843465115
$ws.Cells.Item(13, 15).Value = 1.03
$ws.Cells.Item(13, 16).Value = 1.05186302150144
$ws.Cells.Item(13, 18).Value = 1.045995451256547

# Row 14
$ws.Cells.Item(14, 3).Value = 1.013707209871226
$ws.Cells.Item(14, 4).Value = 1.035561738584025
$ws.Cells.Item(14, 5).Value = 1.029466770581823
$ws.Cells.Item(14, 6).Value = 1.039829308255235
$ws.Cells.Item(14, 9).Value = 1.056841058031452
$ws.Cells.Item(14, 10).Value = 1.042555521440429
$ws.Cells.Item(14, 11).Value = 1.04985892317857
$ws.Cells.Item(14, 12).Value = 1.043869858718199
$ws.Cells.Item(14, 13).Value = 1.054053295897563
$ws.Cells.Item(14, 14).Value = 1.017789123144027
$ws.Cells.Item(14, 15).Value = 1.03
$ws.Cells.Item(14, 16).Value = 1.052599002565867
$ws.Cells.Item(14, 18).Value = 1.046392286943839

# Row 15
$ws.Cells.Item(15, 3).Value = 1.014211477896213
$ws.Cells.Item(15, 4).Value = 1.035917420985465
$ws.Cells.Item(15, 5).Value = 1.029870722634086
$ws.Cells.Item(15, 6).Value = 1.040227354204122
$ws.Cells.Item(15, 9).Value = 1.057014767925865
$ws.Cells.Item(15, 10).Value = 1.042857176965425
$ws.Cells.Item(15, 11).Value = 1.05012400314573
$ws.Cells.Item(15, 12).Value = 1.044181555908167
$ws.Cells.Item(15, 13).Value = 1.054360583563887
$ws.Cells.Item(15, 14).Value = 1.017904979846573
$ws.Cells.Item(15, 15).Value = 1.03
$ws.Cells.Item(15, 16).Value = 1.052879082603679
$ws.Cells.Item(15, 18).Value = 1.046585486070362

# Row 16
$ws.Cells.Item(16, 3).Value = 1.016798521514382
$ws.Cells.Item(16, 4).Value = 1.03765931052541
$ws.Cells.Item(16, 5).Value = 1.031828328351429
$ws.Cells.Item(16, 6).Value = 1.042035421007668
$ws.Cells.Item(16, 9).Value = 1.057816917167687
$ws.Cells.Item(16, 10).Value = 1.044286438977223
$ws.Cells.Item(16, 11).Value = 1.051345568590314
$ws.Cells.Item(16, 12).Value = 1.045610684334922
$ws.Cells.Item(16, 13).Value = 1.055650465291472
$ws.Cells.Item(16, 14).Value = 1.018386137802668
$ws.Cells.Item(16, 15).Value = 1.03
$ws.Cells.Item(16, 16).Value = 1.053860258641476
$ws.Cells.Item(16, 18).Value = 1.047452269554517

# Row 17
$ws.Cells.Item(17, 3).Value = 1.018283055191529
$ws.Cells.Item(17, 4).Value = 1.038625862807284
$ws.Cells.Item(17, 5).Value = 1.032902139959697
$ws.Cells.Item(17, 6).Value = 1.042966281584933
$ws.Cells.Item(17, 9).Value = 1.058238458173882
$ws.Cells.Item(17, 10).Value = 1.045056330411272
$ws.Cells.Item(17, 11).Value = 1.051988918378193
$ws.Cells.Item(17, 12).Value = 1.046356862731122
$ws.Cells.Item(17, 13).Value = 1.056260682424415
$ws.Cells.Item(17, 14).Value = 1.01861585144033
$ws.Cells.Item(17, 15).Value = 1.03
$ws.Cells.Item(17, 16).Value = 1.054214860627538
$ws.Cells.Item(17, 18).Value = 1.047909676267473

# Row 18
$ws.Cells.Item(18, 3).Value = 1.01893822512894
$ws.Cells.Item(18, 4).Value = 1.038987468339449
$ws.Cells.Item(18, 5).Value = 1.033283376278661
$ws.Cells.Item(18, 6).Value = 1.043181849727491
$ws.Cells.Item(18, 9).Value = 1.058356029990853
$ws.Cells.Item(18, 10).Value = 1.045304379053057
$ws.Cells.Item(18, 11).Value = 1.052165530860145
$ws.Cells.Item(18, 12).Value = 1.046551304816285
$ws.Cells.Item(18, 13).Value = 1.056294632180445
$ws.Cells.Item(18, 14).Value = 1.018636193674118
$ws.Cells.Item(18, 15).Value = 1.03
$ws.Cells.Item(18, 16).Value = 1.054006735913576
$ws.Cells.Item(18, 18).Value = 1.048023086283173

# Row 19
$ws.Cells.Item(19, 3).Value = 1.018856639856852
$ws.Cells.Item(19, 4).Value = 1.038818763102341
$ws.Cells.Item(19, 5).Value = 1.033049085760752
$ws.Cells.Item(19, 6).Value = 1.042759062198403
$ws.Cells.Item(19, 9).Value = 1.058210580179989
$ws.Cells.Item(19, 10).Value = 1.045095138428415
$ws.Cells.Item(19, 11).Value = 1.051938328819094
$ws.Cells.Item(19, 12).Value = 1.046259112798651
$ws.Cells.Item(19, 13).Value = 1.055817583460823
$ws.Cells.Item(19, 14).Value = 1.018480245003051
$ws.Cells.Item(19, 15).Value = 1.03
$ws.Cells.Item(19, 16).Value = 1.053308262227438
$ws.Cells.Item(19, 18).Value = 1.047868749071716

# Row 20
$ws.Cells.Item(20, 3).Value = 1.016694643815337
$ws.Cells.Item(20, 4).Value = 1.037148497613182
$ws.Cells.Item(20, 5).Value = 1.031088496033619
$ws.Cells.Item(20, 6).Value = 1.040576564427305
$ws.Cells.Item(20, 9).Value = 1.05732268664698
$ws.Cells.Item(20, 10).Value = 1.043596844873328
$ws.Cells.Item(20, 11).Value = 1.050568140941172
$ws.Cells.Item(20, 12).Value = 1.044605865053911
$ws.Cells.Item(20, 13).Value = 1.053941571282626
$ws.Cells.Item(20, 14).Value = 1.017823628554017
$ws.Cells.Item(20, 15).Value = 1.03
$ws.Cells.Item(20, 16).Value = 1.051303707846876
$ws.Cells.Item(20, 18).Value = 1.046903867618036

# Row 21
$ws.Cells.Item(21, 3).Value = 1.01179745522843
$ws.Cells.Item(21, 4).Value = 1.033798644886887
$ws.Cells.Item(21, 5).Value = 1.027298817750338
$ws.Cells.Item(21, 6).Value = 1.037005822402206
$ws.Cells.Item(21, 9).Value = 1.055747957334901
$ws.Cells.Item(21, 10).Value = 1.040810078565926
$ws.Cells.Item(21, 11).Value = 1.04816678788625
$ws.Cells.Item(21, 12).Value = 1.041780827201398
$ws.Cells.Item(21, 13).Value = 1.051318460959038
$ws.Cells.Item(21, 14).Value = 1.016834295260588
$ws.Cells.Item(21, 15).Value = 1.03
$ws.Cells.Item(21, 16).Value = 1.049187458561952
$ws.Cells.Item(21, 18).Value = 1.045209254675363

# Row 22
$ws.Cells.Item(22, 3).Value = 1.008678198016838
$ws.Cells.Item(22, 4).Value = 1.031674949910564
$ws.Cells.Item(22, 5).Value = 1.024905303363067
$ws.Cells.Item(22, 6).Value = 1.034770882648574
$ws.Cells.Item(22, 9).Value = 1.054741751059238
$ws.Cells.Item(22, 10).Value = 1.039042165946368
$ws.Cells.Item(22, 11).Value = 1.046643994683734
$ws.Cells.Item(22, 12).Value = 1.039999016631155
$ws.Cells.Item(22, 13).Value = 1.049683620282176
$ws.Cells.Item(22, 14).Value = 1.016209572620878
$ws.Cells.Item(22, 15).Value = 1.03
$ws.Cells.Item(22, 16).Value = 1.047893606002904
$ws.Cells.Item(22, 18).Value = 1.044119283205366

# Row 23
$ws.Cells.Item(23, 3).Value = 1.010324315617117
$ws.Cells.Item(23, 4).Value = 1.032789561640935
$ws.Cells.Item(23, 5).Value = 1.026165026208603
$ws.Cells.Item(23, 6).Value = 1.035947559456427
$ws.Cells.Item(23, 9).Value = 1.055269079238716
$ws.Cells.Item(23, 10).Value = 1.039969995993024
$ws.Cells.Item(23, 11).Value = 1.047439843173861
$ws.Cells.Item(23, 12).Value = 1.040934141533466
$ws.Cells.Item(23, 13).Value = 1.050541885109448
$ws.Cells.Item(23, 14).Value = 1.01653408625625
$ws.Cells.Item(23, 15).Value = 1.03
$ws.Cells.Item(23, 16).Value = 1.048572857882689
$ws.Cells.Item(23, 18).Value = 1.044672461676245

# Row 24
$ws.Cells.Item(24, 3).Value = 1.016697017536417
$ws.Cells.Item(24, 4).Value = 1.037125291708805
$ws.Cells.Item(24, 5).Value = 1.031062550448318
$ws.Cells.Item(24, 6).Value = 1.040525488492249
$ws.Cells.Item(24, 9).Value = 1.057299492122203
$ws.Cells.Item(24, 10).Value = 1.043567140681308
$ws.Cells.Item(24, 11).Value = 1.050530371544036
$ws.Cells.Item(24, 12).Value = 1.044565281288241
$ws.Cells.Item(24, 13).Value = 1.05387644018683
$ws.Cells.Item(24, 14).Value = 1.017797764383022
$ws.Cells.Item(24, 15).Value = 1.03
$ws.Cells.Item(24, 16).Value = 1.051211900846272
$ws.Cells.Item(24, 18).Value = 1.046850158109552

# Row 25
$ws.Cells.Item(25, 3).Value = 1.02383866110996
$ws.Cells.Item(25, 4).Value = 1.042005607634499
$ws.Cells.Item(25, 5).Value = 1.036582345777044
$ws.Cells.Item(25, 6).Value = 1.045693669066119
$ws.Cells.Item(25, 9).Value = 1.059535524089374
$ws.Cells.Item(25, 10).Value = 1.04758955653991
$ws.Cells.Item(25, 11).Value = 1.053985047819412
$ws.Cells.Item(25, 12).Value = 1.048638121588143
$ws.Cells.Item(25, 13).Value = 1.0576217681197
$ws.Cells.Item(25, 14).Value = 1.01921033958315
$ws.Cells.Item(25, 15).Value = 1.03
$ws.Cells.Item(25, 16).Value = 1.054176027740823
$ws.Cells.Item(25, 18).Value = 1.049289900085868
